$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row (row 11): correct Right/Wrong marking values
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# "Total" row (row 12): correct total marks and the Max display string
$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "42 / 112"
